$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 'African Trees for Climate Resilience'
$ws.Range("H3").Value = 'SDG 15'
$ws.Range("J3").Value = 'Prof. Guy F. Midgley University of Stellenbosch (gfmidgley@sun.ac.za)'
$ws.Range("P3").Value = 'Powered by:  Stellenbosch University
Catalyzed by: Lacuna-Fund / Meridian (Climate-call) &  FAIR Forward - AI for All
Financed by: BMZ'
$ws.Range("H7").Value = 'SDG 10, SDG 2'
$ws.Range("J7").Value = 'Fundación Despacio (mafe@despacio.org), World Resources Institute'
$ws.Range("O7").Value = 'CC-BY 4.0'
$ws.Range("H8").Value = 'SDG 15'
$ws.Range("J8").Value = 'CTTC María Cuevas (mcuevas@cttc.es),  INVEMAR Cristian Montes (cristian.montes@invemar.org.co)'
$ws.Range("H10").Value = 'SDG 13, SDG 15'
$ws.Range("P10").Value = 'Powered by: Data354
Catalyzed by: FAIR Forward - AI for All (https://www.bmz-digital.global/en/overview-of-initiatives/fair-forward/)
Financed by: BMZ (https://www.bmz-digital.global/en/digital-transformation-and-development-cooperation/)'
$ws.Range("H11").Value = 'SDG 7, SDG 11'
$ws.Range("J11").Value = 'ESPOL University (jecordov@espol.edu.ec)'
$ws.Range("H12").Value = 'SDG 13, SDG 15'
$ws.Range("J12").Value = 'Fundacion Ecociencia (carmenjosse@ecociencia.org)'
$ws.Range("H13").Value = ''
$ws.Range("I13").Value = 'Voice'
$ws.Range("J13").Value = 'Space4Innovation, Diana Mastracci (diana@space4innovation.com)'
$ws.Range("L13").Value = 'Ltome-Katip Indigenous Bioacoustic Dataset
Regions: Samburu (Kenya) · Shuar (Ecuadorian Amazon)
Custodians: Chief Titus Letaapo (Samburu tribe) (Namunyak Conservancy), Chief Mario Vargas Shakaim (Shuar Nation) (MUSAP Biological Station), and Space4Innovation
This dataset contains Indigenous-labelled bioacoustic recordings from two ecosystems—semi-arid savannah and tropical rainforest—collected through AudioMoth bioacustic sensors. Data include species-specific sounds (e.g., elephants, rodents), environmental background, and associated metadata following the CARE Principles for Indigenous Data Governance.
Use cases: biodiversity monitoring, species classification, human–wildlife conflict alerts, and AI model training for conservation.
Limitations: class imbalance (key species overrepresented), environmental noise, and spatial clustering; users should apply noise filtering and ethical review before reuse. These audio data are collected 24/7 when deployed in time periods ranging from hours to several weeks. The data are acquired from multiple microphones spread across the study site. Each microphone has a unique serial number and the geographic locations are provided using GPS. The data are time stamped, however there are data gaps in time and space due to logistics, equipment failure, or power loss. The original data are stored as 16-bit WAV files and are available. To make the data more widely available, they have been uploaded to the Arbimon.org platform. The Arbimon cloud platform is built for bioacoustics analysis using various ML . 

'
$ws.Range("J14").Value = 'Dennis Asamoah Owusu (dowusu@ashesi.edu.gh)'
$ws.Range("H15").Value = 'SDG 15'
$ws.Range("J15").Value = 'Center for Remote Sensing and Geographic Information Services CERSGIS (fkmawusi@gmail.com)'
$ws.Range("H17").Value = 'SDG 2'
$ws.Range("P17").Value = 'Powered by:  RAIL - KNUST (https://rail.knust.edu.gh/)
Catalyzed by: FAIR Forward - AI for All (https://www.bmz-digital.global/en/overview-of-initiatives/fair-forward/)
Financed by: BMZ (https://www.bmz-digital.global/en/digital-transformation-and-development-cooperation/)'
$ws.Range("H18").Value = 'SDG 10'
$ws.Range("H19").Value = 'SDG 2'
$ws.Range("H20").Value = 'SDG 15'
$ws.Range("H21").Value = 'SDG 2, SDG 13'
$ws.Range("H27").Value = 'SDG 15, SDG 13'
$ws.Range("J27").Value = 'Vertify.earth - Michael Anthony (michael@vertify.earth), 
Alsisar Impact - Saurabh Singhavi (saurabh@alsisarimpact.com)'
$ws.Range("H28").Value = 'SDG 2, SDG 10, SDG 5'
$ws.Range("H29").Value = 'SDG 2, SDG 10'
$ws.Range("H30").Value = 'SDG 2'
$ws.Range("H31").Value = 'SDG 10'
$ws.Range("H32").Value = 'SDG 10'
$ws.Range("H33").Value = 'SDG 2'
$ws.Range("H34").Value = 'SDG 2'
$ws.Range("H41").Value = 'SDG 10'
$ws.Range("P41").Value = 'Powered by: Bappenas, Prosa AI
Catalyzed by: FAIR Forward - AI for All
Financed by: BMZ'
$ws.Range("H42").Value = 'SDG 13, SDG 15'
$ws.Range("P42").Value = 'Powered by:  Bappenas
Catalyzed by: FAIR Forward - AI for All
Financed by: BMZ'
$ws.Range("H43").Value = 'SDG 2, SDG 10'
$ws.Range("J43").Value = 'Prosa AI (https://prosa.ai/)'
$ws.Range("H47").Value = 'SDG 2'
$ws.Range("J47").Value = 'International Center of Insect Physiology and Ecology ICIPE (dg@icipe.org)'
$ws.Range("P47").Value = 'Powered by:  International Center of Insect Physiology and Ecology (ICIPE)
Catalyzed by: Lacuna-Fund / Meridian (Climate-call) &  FAIR Forward - AI for All
Financed by: BMZ'
$ws.Range("H53").Value = 'SDG 2, SDG 10'
$ws.Range("P53").Value = 'Powered by: Digital Green, Viamo and Partners, Tech for Her, DynAg, Digifarm, Opportunity International
Catalyzed by: FAIR Forward - AI for All
Financed by: Gates Foundation'
$ws.Range("H54").Value = 'SDG 2, SDG 10'
$ws.Range("P55").Value = 'Powered by:  Association Maidi
Catalyzed by: Lacuna-Fund / Meridian & FAIR Forward - AI for All
Financed by: BMZ'
$ws.Range("H56").Value = 'SDG 10'
$ws.Range("H57").Value = 'SDG 13, SDG 7'
$ws.Range("J57").Value = 'Lahore University of Management Sciences LUMS (naveedarshad@gmail.com)'
$ws.Range("H62").Value = 'SDG 2, SDG 10'
$ws.Range("J62").Value = 'Digital Umuganda (https://digitalumuganda.com/#contact)'
$ws.Range("H63").Value = 'SDG 2'
$ws.Range("J63").Value = 'Benson Kenduiyvo (b.kenduiywo@cgiar.org)'
$ws.Range("J65").Value = 'Digital Umuganda (https://digitalumuganda.com/#contact), FAIR Forward Rwanda'
$ws.Range("P65").Value = 'Powered by: Atingi, Digital Umuganda, Clear Global
Catalyzed by: FAIR Forward - AI for All
Financed by: BMZ'
$ws.Range("H66").Value = 'SDG 2, SDG 13, SDG 11'
$ws.Range("H67").Value = 'SDG 2, SDG 11'
$ws.Range("H68").Value = 'SDG 2'
$ws.Range("H69").Value = 'SDG 2, SDG 11'
$ws.Range("H70").Value = 'SDG 2, SDG 11'
$ws.Range("H71").Value = 'SDG 2, SDG 11, SDG 12'
$ws.Range("H72").Value = 'SDG 2, SDG 11, SDG 13'
$ws.Range("H73").Value = 'SDG 2, SDG 11, SDG 13'
$ws.Range("H74").Value = 'SDG 2, SDG 11, SDG 13'
$ws.Range("H75").Value = 'SDG 2, SDG 13'
$ws.Range("H76").Value = 'SDG 2'
$ws.Range("H77").Value = 'SDG 2, SDG 10'
$ws.Range("H78").Value = 'SDG 2, SDG 13'
$ws.Range("H80").Value = 'SDG 2, SDG 13'
$ws.Range("P80").Value = 'Powered by:  Croppie
Catalyzed by: FAIR Forward - AI for All
Financed by: BMZ'
$ws.Range("J81").Value = 'Sunbird.ai (emwebaze@sunbird.ai)'
$ws.Range("P81").Value = 'Powered by:  Sunbird.ai
Catalyzed by: FAIR Forward - AI for All
Financed by: BMZ'
$ws.Range("H85").Value = 'SDG 13, SDG 5, SDG 10'
